$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 440, pushing existing rows 440-464 down to 441-465.
$ws.Rows("440").Insert()

# Populate the newly inserted row 440 with the new record's data.
$ws.Cells.Item(440, 1).Value = 4
$ws.Cells.Item(440, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(440, 3).Value = "Los Lagos"
$ws.Cells.Item(440, 4).Value = 44931
$ws.Cells.Item(440, 5).Value = 10
$ws.Cells.Item(440, 6).Value = 100112008
$ws.Cells.Item(440, 7).Value = "Coliflor"
$ws.Cells.Item(440, 8).Value = "Sin especificar"
$ws.Cells.Item(440, 9).Value = "Primera"
$ws.Cells.Item(440, 10).Value = 500
$ws.Cells.Item(440, 11).Value = 1500
$ws.Cells.Item(440, 12).Value = 1700
$ws.Cells.Item(440, 13).Value = 1600
$ws.Cells.Item(440, 14).Value = "`$/unidad"
$ws.Cells.Item(440, 15).Value = "Región Metropolitana"
$ws.Cells.Item(440, 16).Value = 1600
$ws.Cells.Item(440, 17).Value = 1
$ws.Cells.Item(440, 18).Value = "Hortaliza"
